# Apply updated dSF (column F) values per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = 0
    11 = 0
    21 = -3
    27 = 0
    29 = 0
    31 = -1
    32 = 2
    34 = -1
    35 = 4
    44 = 0
    50 = -6
    54 = 3
    63 = -2
    65 = -9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
